# Append three new alarm rows (rows 9-11) to the sheet, matching the shape
# of the existing rows (Date, Participant Number, Block Name, Alarm number,
# Vital Sign, Value, Timestamp).
#
# Columns A and B hold values that *look* numeric/date-like ("2024-10-29",
# "10001") but, like the rest of the sheet, must stay plain text. A leading
# apostrophe is the standard Excel way to force text entry and stops
# "2024-10-29" being auto-converted to a date serial number and "10001"
# being auto-converted to a plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "'2024-10-29"
$ws.Range("B9").Value = "'10001"
$ws.Range("C9").Value = "Khushiremote"
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "Heart Rate"
$ws.Range("F9").Value = "High"
$ws.Range("G9").Value = "2024-10-29T14:48:12.956"

# Row 10
$ws.Range("A10").Value = "'2024-10-29"
$ws.Range("B10").Value = "'10001"
$ws.Range("C10").Value = "Khushiremote"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "Blood Pressure"
$ws.Range("F10").Value = "Low"
$ws.Range("G10").Value = "2024-10-29T14:48:13.605"

# Row 11
$ws.Range("A11").Value = "'2024-10-29"
$ws.Range("B11").Value = "'10001"
$ws.Range("C11").Value = "Khushiremote"
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = "Oxygen Saturation"
$ws.Range("F11").Value = "Very High"
$ws.Range("G11").Value = "2024-10-29T14:48:14.571"
